$d = $word.ActiveDocument

# Target cell paragraph originally reads:  "El " + "visualizara la factura"
# across two runs. Final text must read:
#   "El" | " dueño" | " " | "visualizara la factura" | " antes de generarla para su corrección."
# (5 runs, matching the authored diff).
#
# Strategy: do all text insertions first (Find.Execute relocates each time so we
# never rely on stale offsets), then -- once the final text is in place -- force
# the newly-inserted spans to live in their own runs by toggling a character
# property on/off (Word/this engine only coalesces adjacent runs that end up with
# identical *effective* formatting as a side effect of the very write that touches
# them; flipping Bold on then back off forces a genuine, persistent split).

$wdFindContinue = 1

# 1) "El " -> "El"  (drop the trailing space from the first run)
$rng = $d.Content
[void]$rng.Find.Execute("El visualizara la factura", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$start = $rng.Start
$r1 = $d.Range($start, $start + 3)
$r1.Text = "El"

# 2) insert " dueño" right after "El"
$rng = $d.Content
[void]$rng.Find.Execute("Elvisualizara la factura", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$start = $rng.Start
$dueno = " dueño"
$ins = $d.Range($start + 2, $start + 2)
$ins.InsertAfter($dueno)

# 3) insert a single space right after " dueño" (before "visualizara")
$rng = $d.Content
[void]$rng.Find.Execute("El dueñovisualizara la factura", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$start = $rng.Start
$spacePos = $start + 2 + $dueno.Length
$ins2 = $d.Range($spacePos, $spacePos)
$ins2.InsertAfter(" ")

# 4) insert " antes de generarla para su corrección." right after "visualizara la factura"
$tail = " antes de generarla para su corrección."
$rng = $d.Content
[void]$rng.Find.Execute("visualizara la factura", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$vizEnd = $rng.End
$ins3 = $d.Range($vizEnd, $vizEnd)
$ins3.InsertAfter($tail)

# Sanity-check the fully assembled sentence before splitting it into runs.
$full = "El dueño visualizara la factura antes de generarla para su corrección."
$rng = $d.Content
[void]$rng.Find.Execute($full, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $rng.Found) {
    throw "post-insert text does not match expected sentence"
}
$sentenceStart = $rng.Start

# Force-split each newly inserted span into its own run.
function Split-Range($r) {
    $r.Bold = 1
    $r.Bold = 0
}

$duenoStart = $sentenceStart + 2
Split-Range ($d.Range($duenoStart, $duenoStart + $dueno.Length))

$spaceStart = $duenoStart + $dueno.Length
Split-Range ($d.Range($spaceStart, $spaceStart + 1))

$rng = $d.Content
[void]$rng.Find.Execute("visualizara la factura", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$tailStart = $rng.End
Split-Range ($d.Range($tailStart, $tailStart + $tail.Length))

Write-Output "done"
